# The sheet holds a weekly price log for "Cebollín baby" (Agrícola del
# Norte S.A. de Arica). A new weekly observation needs to be inserted as
# row 18, pushing the existing rows 18-51 down to 19-52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 18 (shifts 18..51 -> 19..52).
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(18, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(18, 4).Value = 44487
$ws.Cells.Item(18, 5).Value = 15
$ws.Cells.Item(18, 6).Value = 100112038
$ws.Cells.Item(18, 7).Value = "Cebollín baby"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 270
$ws.Cells.Item(18, 11).Value = 1000
$ws.Cells.Item(18, 12).Value = 1200
$ws.Cells.Item(18, 13).Value = 1100
$ws.Cells.Item(18, 14).Value = "`$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(18, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(18, 16).Value = 550
$ws.Cells.Item(18, 17).Value = 2
$ws.Cells.Item(18, 18).Value = "Hortaliza"
